# Adapt column header formatting to respective input file names.
# "_old" -> "_FV2310" and "_new" -> "_FV2404" suffixes on the first 20
# header cells (row 1), then turn the used range into a real Excel table
# (Table1) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2310 = "_FV2310"
$fv2404 = "_FV2404"

# Row 1 headers, columns A..U (1..21). Column K ("diff") has no suffix.
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = [string]$cell.Value2

    if ($text.EndsWith($oldSuffix)) {
        $cell.Value = $text.Substring(0, $text.Length - $oldSuffix.Length) + $fv2310
    } elseif ($text.EndsWith($newSuffix)) {
        $cell.Value = $text.Substring(0, $text.Length - $newSuffix.Length) + $fv2404
    }
}

# Convert the data range into a native Excel Table ("Table1") so the
# header row doubles as column metadata / autofilter.
$listRange = $ws.Range("A1:U64")
$table = $ws.ListObjects.Add(1, $listRange, $null, 1)
$table.Name = "Table1"

# Freeze the header row (row 2 onward scrolls, row 1 stays pinned).
$ws.Activate()
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.SplitColumn = 0
$excel.ActiveWindow.FreezePanes = $true
